$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that is refreshed on every
# export. Update it for every data row (row 2 through the last used row).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 452 }

$ws.Range("C2:C$lastRow").Value = 45175
